# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt - Coliflor".
# The new record is inserted as row 365 (pushing the former rows 365..445 down to
# 366..446), matching the source diff which adds one new dated observation and
# shifts every subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 365..445 down to 366..446, leaving a fresh blank row 365.
$ws.Rows.Item(365).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(365, 1).Value  = 4
$ws.Cells.Item(365, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(365, 3).Value  = "Los Lagos"
$ws.Cells.Item(365, 4).Value  = 44889
$ws.Cells.Item(365, 5).Value  = 10
$ws.Cells.Item(365, 6).Value  = 100112008
$ws.Cells.Item(365, 7).Value  = "Coliflor"
$ws.Cells.Item(365, 8).Value  = "Sin especificar"
$ws.Cells.Item(365, 9).Value  = "Primera"
$ws.Cells.Item(365, 10).Value = 600
$ws.Cells.Item(365, 11).Value = 1500
$ws.Cells.Item(365, 12).Value = 1500
$ws.Cells.Item(365, 13).Value = 1500
$ws.Cells.Item(365, 14).Value = "`$/unidad"
$ws.Cells.Item(365, 15).Value = "Región Metropolitana"
$ws.Cells.Item(365, 16).Value = 1500
$ws.Cells.Item(365, 17).Value = 1
$ws.Cells.Item(365, 18).Value = "Hortaliza"
